# Add 2025-01-28 year-to-date violent crime counts (column L = year 2025).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 12).Value = 394
$ws.Cells.Item(3, 12).Value = 421
$ws.Cells.Item(4, 12).Value = 110
$ws.Cells.Item(6, 12).Value = 496
$ws.Cells.Item(7, 12).Value = 1449

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(4, 12).Value = 5
$ws.Cells.Item(8, 12).Value = 78
$ws.Cells.Item(12, 12).Value = 4
$ws.Cells.Item(20, 12).Value = 37
$ws.Cells.Item(27, 12).Value = 11
$ws.Cells.Item(29, 12).Value = 81
$ws.Cells.Item(31, 12).Value = 19
$ws.Cells.Item(34, 12).Value = 10
$ws.Cells.Item(42, 12).Value = 50
$ws.Cells.Item(43, 12).Value = 14
$ws.Cells.Item(44, 12).Value = 7
$ws.Cells.Item(53, 12).Value = 21
$ws.Cells.Item(54, 12).Value = 26
$ws.Cells.Item(60, 12).Value = 9
$ws.Cells.Item(63, 12).Value = 5
$ws.Cells.Item(65, 12).Value = 29
$ws.Cells.Item(67, 12).Value = 46
$ws.Cells.Item(73, 12).Value = 11
$ws.Cells.Item(78, 12).Value = 15
$ws.Cells.Item(79, 12).Value = 42
$ws.Cells.Item(83, 12).Value = 33
$ws.Cells.Item(85, 12).Value = 77
$ws.Cells.Item(89, 12).Value = 16
$ws.Cells.Item(90, 12).Value = 11
$ws.Cells.Item(91, 12).Value = 19
$ws.Cells.Item(95, 12).Value = 18
$ws.Cells.Item(96, 12).Value = 12
$ws.Cells.Item(97, 12).Value = 21
$ws.Cells.Item(101, 12).Value = 1449

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(2, 12).Value = 7
$ws.Cells.Item(7, 12).Value = 12

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(4, 12).Value = 2
$ws.Cells.Item(7, 12).Value = 16

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 77

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(3, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 21

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(3, 12).Value = 24
$ws.Cells.Item(7, 12).Value = 78

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(3, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 33

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(2, 12).Value = 7
$ws.Cells.Item(7, 12).Value = 18

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 12).Value = 7
$ws.Cells.Item(3, 12).Value = 6
$ws.Cells.Item(6, 12).Value = 13
$ws.Cells.Item(7, 12).Value = 29

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(2, 12).Value = 4
$ws.Cells.Item(7, 12).Value = 19

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(6, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 46

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(3, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 26

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(3, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 81

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 7

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 12).Value = 10
$ws.Cells.Item(3, 12).Value = 13
$ws.Cells.Item(6, 12).Value = 23
$ws.Cells.Item(7, 12).Value = 50

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(2, 12).Value = 4
$ws.Cells.Item(6, 12).Value = 15

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 12).Value = 8
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(4, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 19

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(3, 12).Value = 14
$ws.Cells.Item(6, 12).Value = 10
$ws.Cells.Item(7, 12).Value = 42

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 12).Value = 12
$ws.Cells.Item(6, 12).Value = 12
$ws.Cells.Item(7, 12).Value = 37

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(6, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 10

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(2, 12).Value = 5
$ws.Cells.Item(7, 12).Value = 11

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(6, 12).Value = 19
$ws.Cells.Item(7, 12).Value = 21

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(5, 12).Value = 5
$ws.Cells.Item(6, 12).Value = 11

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(2, 12).Value = 6
$ws.Cells.Item(3, 12).Value = 3
$ws.Cells.Item(7, 12).Value = 11

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(2, 12).Value = 3
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(7, 12).Value = 9

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(2, 12).Value = 4
$ws.Cells.Item(3, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 14

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Cells.Item(5, 12).Value = 3
$ws.Cells.Item(6, 12).Value = 5

$ws = $wb.Worksheets.Item("Beverly")
$ws.Cells.Item(2, 12).Value = 2
$ws.Cells.Item(6, 12).Value = 4
